$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column "K" (header in G1) values for rows 2-8 per the regenerated
# save_data (K replaces the old Strike# values, recalculated from std/mean).
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 1
